$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.071.99'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '2.012.89'
$ws.Range('E3').Value = '  -1.79%  '
$ws.Range('E4').Value = '  -0.16%  '
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.606'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.26'
$ws.Range('E8').Value = '  -1.48%  '
$ws.Range('E9').Value = '  -2.86%  '
$ws.Range('E10').Value = '  -4.18%  '
$ws.Range('E11').Value = '  -4.01%  '
$ws.Range('D12').Value = '2.310.55'
$ws.Range('E12').Value = '  -1.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.05'
$ws.Range('E13').Value = '  -3.11%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '19.76'
$ws.Range('E14').Value = '  -4.20%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.735'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.16'
$ws.Range('E16').Value = '  -1.76%  '
$ws.Range('D17').Value = '2.014.26'
$ws.Range('E17').Value = '  -1.75%  '
$ws.Range('D18').Value = '36.987.05'
$ws.Range('E18').Value = '  -0.83%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.19'
$ws.Range('E19').Value = '  +1.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '68.29'
$ws.Range('E20').Value = '  -2.27%  '
$ws.Range('E21').Value = '  -4.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '223.43'
$ws.Range('E22').Value = '  -1.15%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +1.94%  '
$ws.Range('E25').Value = '  -4.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.48'
$ws.Range('E26').Value = '  -2.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.94'
$ws.Range('E27').Value = '  -6.41%  '
$ws.Range('E28').Value = '  -3.53%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '18.58'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('E30').Value = '  -7.18%  '
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('E32').Value = '  -2.90%  '
$ws.Range('E33').Value = '  -1.96%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.44'
$ws.Range('E34').Value = '  -2.78%  '
$ws.Range('E35').Value = '  -3.02%  '
$ws.Range('E36').Value = '  +1.97%  '
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.14'
$ws.Range('E38').Value = '  -1.47%  '
$ws.Range('E39').Value = '  -1.45%  '
$ws.Range('D40').Value = '1.462.70'
$ws.Range('E40').Value = '  -2.95%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0212'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.29'
$ws.Range('E42').Value = '  +20.89%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '94.52'
$ws.Range('E43').Value = '  -1.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0911'
$ws.Range('E44').Value = '  -2.34%  '
$ws.Range('E45').Value = '  -4.19%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '15.90'
$ws.Range('E46').Value = '  -5.90%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.12'
$ws.Range('E47').Value = '  -2.02%  '
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.06'
$ws.Range('E49').Value = '  -1.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.89'
$ws.Range('E50').Value = '  -0.87%  '
$ws.Range('D51').Value = '2.198.73'
$ws.Range('E51').Value = '  -1.80%  '
